$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C3"   = -11.4162
    "C14"  = -13.60769999999998
    "C16"  = -13.83679999999999
    "C21"  = -12.3496
    "C23"  = -12.6208
    "C25"  = -13.24719999999999
    "C26"  = -12.60320000000001
    "C29"  = -10.97910000000001
    "C40"  = -13.0912
    "C53"  = -10.53300000000001
    "C57"  = -14.0092
    "C59"  = -12.6101
    "C65"  = -12.6008
    "C69"  = -11.6664
    "C79"  = -10.73580000000001
    "C83"  = -13.6211
    "C91"  = -10.2512
    "C93"  = -11.64950000000001
    "C100" = -13.23569999999999
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
